$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (changed) date column C for rows 2-7
# from serial date 45224 (2023-10-25) to 45233 (2023-11-03)
$ws.Range("C2:C7").Value = 45233
